$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "[-, -, -, 'MCT-3A-Eletrohidráulica']"
$ws.Range("E4").Value = "-"
$ws.Range("E6").Value = "-"
$ws.Range("D7").Value = "[-, -, -, 'MCT-3A-Eletrohidráulica']"
$ws.Range("E7").Value = "['MCT-3A-Eletrohidráulica', -, -, -]"
$ws.Range("E8").Value = "['MCT-3A-Eletrohidráulica', -, -, -]"
